$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.567399999999993
$ws.Range("D3").Value = -5.800999999999999
$ws.Range("D5").Value = -7.607699999999995
$ws.Range("E7").Value = 11.9393
$ws.Range("B9").Value = 8.750700000000007
$ws.Range("E9").Value = 14.40930000000002
$ws.Range("D11").Value = -8.4518
$ws.Range("D12").Value = -8.373300000000002
$ws.Range("B13").Value = 5.853499999999998
$ws.Range("B16").Value = 9.155800000000006
$ws.Range("B18").Value = 5.213299999999998
$ws.Range("B20").Value = 5.3765
$ws.Range("D21").Value = -7.774700000000003
$ws.Range("E21").Value = 13.19820000000001
